# Commit: Fri, May 08, 2020 12:07:05 PM
#
# 1) Three tables (slides 14, 15, 16) switch from table style
#    {9474CB1B-F018-4EC4-A4F5-C2B057BFC6B3} ("No Style, No Grid") to
#    {30D4788F-2E10-4929-8818-65557C26E32C}.
# 2) The presentation's theme (ppt/theme/theme1.xml, "Integral" / "Red
#    Violet") is recoloured to the stock "Office" palette (the palette
#    that, before the edit, only lived in the Notes Master's theme
#    part, ppt/theme/theme2.xml). The font scheme and format scheme are
#    already identical between the two theme parts, so only the 12
#    scheme colours need to change.

$p = $ppt.ActivePresentation

$newStyleId = "{30D4788F-2E10-4929-8818-65557C26E32C}"
$tableSlides = @(14, 15, 16)
foreach ($slideIdx in $tableSlides) {
    $slide = $p.Slides.Item($slideIdx)
    for ($shapeIdx = 1; $shapeIdx -le $slide.Shapes.Count; $shapeIdx++) {
        $shape = $slide.Shapes.Item($shapeIdx)
        if ($shape.HasTable) {
            $shape.Table.ApplyStyle($newStyleId)
        }
    }
}

# Recolour the (single, shared) theme to the "Office" colour scheme.
$colorScheme = $p.Slides.Item(1).ThemeColorScheme
$colorScheme.Item(1).RGB = 0          # dk1      000000
$colorScheme.Item(2).RGB = 16777215   # lt1      FFFFFF
$colorScheme.Item(3).RGB = 6968388    # dk2      44546A
$colorScheme.Item(4).RGB = 15132391   # lt2      E7E6E6
$colorScheme.Item(5).RGB = 13998939   # accent1  5B9BD5
$colorScheme.Item(6).RGB = 3243501    # accent2  ED7D31
$colorScheme.Item(7).RGB = 10855845   # accent3  A5A5A5
$colorScheme.Item(8).RGB = 49407      # accent4  FFC000
$colorScheme.Item(9).RGB = 12874308   # accent5  4472C4
$colorScheme.Item(10).RGB = 4697456   # accent6  70AD47
$colorScheme.Item(11).RGB = 12673797  # hlink    0563C1
$colorScheme.Item(12).RGB = 7491477   # folHlink 954F72
